$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$c8 = $ws1.Range("C8")
$c8.Hyperlinks.Add($c8, "", "'Sheet''s 2'!InSheetName", "", "'Sheet''s 2'!InSheetName")
